$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the school (HoiDong) name for both data rows: "THCS Lý Thường Kiệt" -> "THCS Phú Ninh"
$ws.Range("AG2").Value = "THCS Phú Ninh"
$ws.Range("AG3").Value = "THCS Phú Ninh"

# Update the view: scroll so column R is the left-most visible column, and move the
# active cell/selection to AG9
$ws.Application.ActiveWindow.ScrollColumn = 18
$ws.Range("AG9").Select()
